# "Generate Report for Handoff" - refresh the localization-status report
# with a new handoff run: new guid-named source/xliff files and new
# generate/handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldGuid = "56689279-fb5f-423f-a85b-0b22247fe4df"
$newGuid = "48d916f5-a784-43f5-9317-833c35d845f1"

$oldFileName = "$oldGuid.md"
$newFileName = "$newGuid.md"

$oldPath = "e2e\$oldGuid.md"
$newPath = "e2e\$newGuid.md"

$newGenerateDate = "2017-01-03 08:49:44"

$oldZhXlf = "$oldGuid.4c6b2016e0d711cbd266ac9a610710e21d199d55.zh-cn.xlf"
$newZhXlf = "$newGuid.a7b702a63d88791c80fae4f01594805d9a3d2ede.zh-cn.xlf"
$newZhHandoffDate = "2017-01-03 08:49:34"

$oldDeXlf = "$oldGuid.4c6b2016e0d711cbd266ac9a610710e21d199d55.de-de.xlf"
$newDeXlf = "$newGuid.a7b702a63d88791c80fae4f01594805d9a3d2ede.de-de.xlf"

# The external hyperlink target (the commit blob URL) itself is not part
# of this change, only the guid-based display text shown in the cells.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/c2ae70ff7ca8301501f1207d244d0dbd321eb8bb/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Overview sheet: A2 (File Name), B2 (Path And Name, hyperlinked),
# G2 (Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPath
$wsOverview.Range("G2").Value = $newGenerateDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $newPath)

# ---------------------------------------------------------------------
# zh-cn sheet: A2 (Source File Name, hyperlinked), G2 (Latest Handoff
# File), H2 (Latest Handoff Datetime)
# ---------------------------------------------------------------------
$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhHandoffDate

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $newFileName)

# ---------------------------------------------------------------------
# de-de sheet: A2 (Source File Name, hyperlinked), G2 (Latest Handoff
# File), H2 (Latest Handoff Datetime - same value as Overview!G2)
# ---------------------------------------------------------------------
$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $newGenerateDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $newFileName)
